# Insert a new weekly data row at row 64 (pushes existing rows 64..204 down to 65..205)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new week's record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R repeat the same market/category metadata
# as the row below (old row 64, now row 65); D (Fecha) and J (Volumen) are
# the genuinely new values for this entry.
$ws.Cells.Item(64, 1).Value = 5
$ws.Cells.Item(64, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value = "Maule"
$ws.Cells.Item(64, 4).Value = 44519
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(64, 6).Value = 100112009
$ws.Cells.Item(64, 7).Value = "Acelga"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 2000
$ws.Cells.Item(64, 12).Value = 2000
$ws.Cells.Item(64, 13).Value = 2000
$ws.Cells.Item(64, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 500
$ws.Cells.Item(64, 17).Value = 4
$ws.Cells.Item(64, 18).Value = "Hortaliza"
